$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 4 and 5 (entire rows), then update remaining values
$ws.Rows("4:5").Delete()

$ws.Range("A2").Value = 252417
$ws.Range("A3").Value = 252980
